$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows whose balances changed (old positions, highest row number first)
$ws.Rows.Item(246).Delete()   # 001090818 MARIA -2.97
$ws.Rows.Item(234).Delete()   # 004335144 EDMUNDO 0.1
$ws.Rows.Item(103).Delete()   # 004382374 THEOMAR 46.89
$ws.Rows.Item(22).Delete()    # 002823185 SIMONE 98.41

# Insert rows at their new sorted positions with updated balances
$ws.Rows.Item(18).Insert()
$ws.Cells.Item(18,1).Value = "'004382374"
$ws.Cells.Item(18,2).Value = "THEOMAR"
$ws.Cells.Item(18,3).Value = 129.48

$ws.Rows.Item(19).Insert()
$ws.Cells.Item(19,1).Value = "'001090818"
$ws.Cells.Item(19,2).Value = "MARIA"
$ws.Cells.Item(19,3).Value = 105.14

$ws.Rows.Item(20).Insert()
$ws.Cells.Item(20,1).Value = "'002823185"
$ws.Cells.Item(20,2).Value = "SIMONE"
$ws.Cells.Item(20,3).Value = 101.16

$ws.Rows.Item(204).Insert()
$ws.Cells.Item(204,1).Value = "'004335144"
$ws.Cells.Item(204,2).Value = "EDMUNDO"
$ws.Cells.Item(204,3).Value = 2.55

